{"js": "// Helper: wrap a body-fragment of OOXML markup into the \"Flat OPC\" package\n// format required by Range.insertOoxml() / Body.insertOoxml().\nfunction flatOpc(bodyInnerXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>${bodyInnerXml}</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two anchor paragraphs by their exact text.\nlet afterAttributesPara = null;\nlet networkingPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Resource nodes don\\u2019t have premade attributes\") {\n    afterAttributesPara = paragraphs.items[i];\n  } else if (t === \"Very little knowledge on networking\") {\n    networkingPara = paragraphs.items[i];\n  }\n}\nif (!afterAttributesPara) throw new Error(\"Could not find 'Resource nodes...' paragraph\");\nif (!networkingPara) throw new Error(\"Could not find 'Very little knowledge on networking' paragraph\");\n\n// --- 1) Insert four new paragraphs after \"Resource nodes don't have premade attributes\" ---\nconst newParagraphsXml = `\n  <w:p>\n    <w:r><w:t xml:space=\"preserve\">No longer uses </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>enum</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> for resource and collection method identification</w:t></w:r>\n  </w:p>\n  <w:p>\n    <w:r><w:t>Resource, Collection Method, and Attributes are now Objects containing structs with variables</w:t></w:r>\n  </w:p>\n  <w:p>\n    <w:r><w:t xml:space=\"preserve\">Collection Method no longer has </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>resourceUsedFor</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> attribute</w:t></w:r>\n  </w:p>\n  <w:p>\n    <w:r><w:t>Attributes are no longer 4 lists of different variable</w:t></w:r>\n    <w:r><w:t xml:space=\"preserve\">s, instead an object with a struct containing a </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>Varaible</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> object and a string</w:t></w:r>\n  </w:p>`;\n\n// Create one placeholder paragraph right after the anchor, then replace its\n// content (the whole paragraph range) with the four real paragraphs above.\nconst insertionAnchor = afterAttributesPara.insertParagraph(\"\", \"After\");\nawait context.sync();\ninsertionAnchor.getRange().insertOoxml(flatOpc(newParagraphsXml), \"Replace\");\nawait context.sync();\n\n// --- 2) Split \"Very little knowledge on networking\" into three runs, and ---\n// --- 3) add a new \"Overscoped\" paragraph right after it.                 ---\n// Re-query the paragraph collection since the previous insert shifted the\n// document and any previously captured paragraph objects are now stale.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\nnetworkingPara = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"Very little knowledge on networking\") {\n    networkingPara = paragraphs2.items[i];\n    break;\n  }\n}\nif (!networkingPara) throw new Error(\"Could not re-find 'Very little knowledge on networking' paragraph\");\n\nconst networkingAndOverscopedXml = `\n  <w:p>\n    <w:r><w:t>Very little knowledge o</w:t></w:r>\n    <w:r><w:t>f</w:t></w:r>\n    <w:r><w:t xml:space=\"preserve\"> networking</w:t></w:r>\n  </w:p>\n  <w:p>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>Overscoped</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n  </w:p>`;\n\nnetworkingPara.getRange().insertOoxml(flatOpc(networkingAndOverscopedXml), \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) script.\n# Document is already open as $word.ActiveDocument ($d below).\n#\n# The edit:\n#   1) Insert four new paragraphs right after\n#      \"Resource nodes don't have premade attributes\".\n#   2) Split \"Very little knowledge on networking\" into three runs\n#      (\"Very little knowledge o\" + \"f\" + \" networking\").\n#   3) Add a new \"Overscoped\" paragraph right after that one.\n\n$d = $word.ActiveDocument\n\nfunction New-FlatOpc([string]$bodyInnerXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyInnerXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# --- 1) Locate \"Resource nodes don't have premade attributes\" and insert ---\n# --- the four new paragraphs right after it.                            ---\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"Resource nodes don\") | Out-Null\n$attributesPara = $findRange.Paragraphs(1)\n\n# Create a new empty paragraph right after it, then fill that paragraph's\n# range with the four real paragraphs via InsertXML (this replaces the\n# (empty) contents of the range it's called on).\n$attributesPara.Range.InsertParagraphAfter()\n$newParasAnchor = $attributesPara.Next()\n\n$newParagraphsXml = @'\n<w:p>\n  <w:r><w:t xml:space=\"preserve\">No longer uses </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>enum</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> for resource and collection method identification</w:t></w:r>\n</w:p>\n<w:p>\n  <w:r><w:t>Resource, Collection Method, and Attributes are now Objects containing structs with variables</w:t></w:r>\n</w:p>\n<w:p>\n  <w:r><w:t xml:space=\"preserve\">Collection Method no longer has </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>resourceUsedFor</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> attribute</w:t></w:r>\n</w:p>\n<w:p>\n  <w:r><w:t>Attributes are no longer 4 lists of different variable</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">s, instead an object with a struct containing a </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>Varaible</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> object and a string</w:t></w:r>\n</w:p>\n'@\n\n$newParasAnchor.Range.InsertXML((New-FlatOpc $newParagraphsXml))\n\n# --- 2) Split \"Very little knowledge on networking\" into three runs, and ---\n# --- 3) add a new \"Overscoped\" paragraph right after it.                 ---\n$findRange2 = $d.Content\n$findRange2.Find.ClearFormatting()\n$findRange2.Find.Execute(\"Very little knowledge on networking\") | Out-Null\n$networkingPara = $findRange2.Paragraphs(1)\n\n$networkingAndOverscopedXml = @'\n<w:p>\n  <w:r><w:t>Very little knowledge o</w:t></w:r>\n  <w:r><w:t>f</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> networking</w:t></w:r>\n</w:p>\n<w:p>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>Overscoped</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n</w:p>\n'@\n\n$networkingPara.Range.InsertXML((New-FlatOpc $networkingAndOverscopedXml))\n"}
